# correções nos métodos listar
$wb = $excel.ActiveWorkbook

# --- Sheet "pesquisadores" ---
$ws2 = $wb.Worksheets.Item("pesquisadores")

# A3 was stored as text "2"; make it a real number
$ws2.Range("A3").Value = 2

# New row 4: Ulises
$ws2.Range("A4").Value = 2
$ws2.Range("B4").Value = "Ulises"
$ws2.Range("C4").Value = "ulises"
$ws2.Range("D4").Value = "ulises"

# New row 5: Jose (id kept as text "5")
$ws2.Range("A5").Value = "'5"
$ws2.Range("B5").Value = "Jose"
$ws2.Range("C5").Value = "jose"
$ws2.Range("D5").Value = "jose"

# --- Sheet "fisioterapeutas" ---
$ws3 = $wb.Worksheets.Item("fisioterapeutas")

# A2 was stored as text "1"; make it a real number
$ws3.Range("A2").Value = 1

# New row 3: Abraoo
$ws3.Range("A3").Value = 3
$ws3.Range("B3").Value = "Abraoo"
$ws3.Range("C3").Value = "abraoo"
$ws3.Range("D3").Value = "abraoo"

# New row 4: Jose (id kept as text "5")
$ws3.Range("A4").Value = "'5"
$ws3.Range("B4").Value = "Jose"
$ws3.Range("C4").Value = "jose"
$ws3.Range("D4").Value = "jose"
